$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the Price (D) column as Text so numeric-looking strings
# (e.g. "1.002") are stored as text, matching the original inline-string cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '22.404.10'
$ws.Range('E2').Value = '  -4.80%  '

$ws.Range('D3').Value = '1.573.41'
$ws.Range('E3').Value = '  -4.71%  '

$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.16%  '

$ws.Range('D5').Value = '1.002'
$ws.Range('E5').Value = '  +0.14%  '

$ws.Range('D6').Value = '291.14'
$ws.Range('E6').Value = '  -3.10%  '

$ws.Range('D7').Value = '0.3672'
$ws.Range('E7').Value = '  -3.06%  '

$ws.Range('D8').Value = '49.58'
$ws.Range('E8').Value = '  -2.41%  '

$ws.Range('D9').Value = '0.3359'
$ws.Range('E9').Value = '  -5.99%  '

$ws.Range('D10').Value = '1.170'
$ws.Range('E10').Value = '  -4.76%  '

$ws.Range('D11').Value = '0.07582'
$ws.Range('E11').Value = '  -6.56%  '

$ws.Range('D12').Value = '1.002'
$ws.Range('E12').Value = '  +0.13%  '

$ws.Range('D13').Value = '21.08'
$ws.Range('E13').Value = '  -4.86%  '

$ws.Range('D14').Value = '6.080'
$ws.Range('E14').Value = '  -5.43%  '

$ws.Range('D15').Value = '6.877'
$ws.Range('E15').Value = '  -7.56%  '

$ws.Range('D16').Value = '1.571.42'
$ws.Range('E16').Value = '  -4.92%  '

$ws.Range('D17').Value = '0.00001134'
$ws.Range('E17').Value = '  -6.01%  '

$ws.Range('D18').Value = '89.43'
$ws.Range('E18').Value = '  -8.03%  '

$ws.Range('D19').Value = '0.06747'
$ws.Range('E19').Value = '  -3.46%  '

$ws.Range('D20').Value = '1.002'
$ws.Range('E20').Value = '  +0.23%  '

$ws.Range('D21').Value = '6.248'
$ws.Range('E21').Value = '  -8.02%  '

$ws.Range('D22').Value = '16.32'
$ws.Range('E22').Value = '  -6.72%  '

$ws.Range('D23').Value = '11.91'
$ws.Range('E23').Value = '  -5.79%  '

$ws.Range('D24').Value = '22.425.63'
$ws.Range('E24').Value = '  -4.78%  '

$ws.Range('D25').Value = '2.410'
$ws.Range('E25').Value = '  -3.01%  '

$ws.Range('D26').Value = '2.954'
$ws.Range('E26').Value = '  +0.76%  '

$ws.Range('D27').Value = '19.76'
$ws.Range('E27').Value = '  -6.08%  '

$ws.Range('D28').Value = '145.44'
$ws.Range('E28').Value = '  -4.63%  '

$ws.Range('D29').Value = '4.942'
$ws.Range('E29').Value = '  -5.62%  '

$ws.Range('D30').Value = '124.69'
$ws.Range('E30').Value = '  -6.53%  '

$ws.Range('D31').Value = '1.747.55'
$ws.Range('E31').Value = '  -4.61%  '

$ws.Range('D32').Value = '6.255'
$ws.Range('E32').Value = '  -10.90%  '

$ws.Range('D33').Value = '1.976'
$ws.Range('E33').Value = '  -7.73%  '

$ws.Range('D34').Value = '0.9708'
$ws.Range('E34').Value = '  -6.28%  '

$ws.Range('D35').Value = '10.44'
$ws.Range('E35').Value = '  -12.25%  '

$ws.Range('D36').Value = '0.08481'
$ws.Range('E36').Value = '  -2.65%  '

$ws.Range('D37').Value = '0.02523'
$ws.Range('E37').Value = '  -8.17%  '

$ws.Range('D38').Value = '0.2290'
$ws.Range('E38').Value = '  -6.74%  '

$ws.Range('D39').Value = '0.06527'
$ws.Range('E39').Value = '  -5.38%  '

$ws.Range('D40').Value = '5.473'
$ws.Range('E40').Value = '  -8.99%  '

$ws.Range('B41').Value = 'Aptos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D41').Value = '11.82'
$ws.Range('E41').Value = '  -10.58%  '

$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = '1.264'
$ws.Range('E42').Value = '  -4.55%  '

$ws.Range('D43').Value = '0.6363'
$ws.Range('E43').Value = '  -8.29%  '

$ws.Range('E44').Value = '  -7.62%  '

$ws.Range('D45').Value = '1.001'
$ws.Range('E45').Value = '  +0.12%  '

$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = '0.5985'
$ws.Range('E46').Value = '  -7.34%  '

$ws.Range('B47').Value = 'PancakeSwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D47').Value = '3.775'
$ws.Range('E47').Value = '  -3.98%  '

$ws.Range('D48').Value = '2.118'
$ws.Range('E48').Value = '  -6.89%  '

$ws.Range('D49').Value = '120.85'
$ws.Range('E49').Value = '  -5.95%  '

$ws.Range('D50').Value = '0.07262'
$ws.Range('E50').Value = '  -7.18%  '

$ws.Range('D51').Value = '1.187'
$ws.Range('E51').Value = '  -0.67%  '

# Reset the Price column style back to Normal so no extra formatting
# is attached to the cells (keeps the cells looking like plain text cells).
$ws.Range("D2:D51").Style = "Normal"